$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cells in existing rows 375-422 (new week inserted, data shifted) ---
# Row 375
$ws.Range("D375").Value = 44522
$ws.Range("K375").Value = 12000
$ws.Range("L375").Value = 13000
$ws.Range("M375").Value = 12500
$ws.Range("P375").Value = 833

# Row 376
$ws.Range("D376").Value = 44522
$ws.Range("J376").Value = 120
$ws.Range("K376").Value = 10000
$ws.Range("L376").Value = 11000
$ws.Range("M376").Value = 10500
$ws.Range("P376").Value = 700

# Row 377
$ws.Range("D377").Value = 44522
$ws.Range("J377").Value = 130
$ws.Range("K377").Value = 8000
$ws.Range("L377").Value = 9000
$ws.Range("M377").Value = 8500
$ws.Range("P377").Value = 567

# Row 378
$ws.Range("H378").Value = "Zafiro rojo"
$ws.Range("K378").Value = 37000
$ws.Range("L378").Value = 38000
$ws.Range("M378").Value = 37500
$ws.Range("P378").Value = 2500

# Row 379
$ws.Range("H379").Value = "Zafiro rojo"
$ws.Range("K379").Value = 34000
$ws.Range("L379").Value = 35000
$ws.Range("M379").Value = 34500
$ws.Range("P379").Value = 2300

# Row 380
$ws.Range("H380").Value = "Zafiro rojo"
$ws.Range("K380").Value = 28000
$ws.Range("L380").Value = 30000
$ws.Range("M380").Value = 29000
$ws.Range("P380").Value = 1933

# Row 381
$ws.Range("D381").Value = 44491
$ws.Range("H381").Value = "Zafiro verde"
$ws.Range("K381").Value = 27000
$ws.Range("L381").Value = 28000
$ws.Range("M381").Value = 27500
$ws.Range("P381").Value = 1833

# Row 382
$ws.Range("D382").Value = 44491
$ws.Range("H382").Value = "Zafiro verde"
$ws.Range("K382").Value = 24000
$ws.Range("L382").Value = 25000
$ws.Range("M382").Value = 24500
$ws.Range("P382").Value = 1633

# Row 383
$ws.Range("D383").Value = 44491
$ws.Range("H383").Value = "Zafiro verde"
$ws.Range("K383").Value = 22000
$ws.Range("L383").Value = 23000
$ws.Range("M383").Value = 22500
$ws.Range("P383").Value = 1500

# Row 384
$ws.Range("D384").Value = 44225
$ws.Range("K384").Value = 9000
$ws.Range("L384").Value = 10000
$ws.Range("M384").Value = 9500
$ws.Range("P384").Value = 633

# Row 385
$ws.Range("D385").Value = 44225
$ws.Range("K385").Value = 7000
$ws.Range("L385").Value = 8000
$ws.Range("M385").Value = 7500
$ws.Range("P385").Value = 500

# Row 386
$ws.Range("D386").Value = 44225
$ws.Range("J386").Value = 160
$ws.Range("K386").Value = 5500
$ws.Range("L386").Value = 6000
$ws.Range("M386").Value = 5750
$ws.Range("P386").Value = 383

# Row 387
$ws.Range("D387").Value = 44232
$ws.Range("K387").Value = 11000
$ws.Range("L387").Value = 12000
$ws.Range("M387").Value = 11500
$ws.Range("P387").Value = 767

# Row 388
$ws.Range("D388").Value = 44232
$ws.Range("J388").Value = 140
$ws.Range("K388").Value = 9000
$ws.Range("L388").Value = 10000
$ws.Range("M388").Value = 9571
$ws.Range("P388").Value = 638

# Row 389
$ws.Range("D389").Value = 44232
$ws.Range("J389").Value = 130
$ws.Range("K389").Value = 7000
$ws.Range("L389").Value = 8000
$ws.Range("M389").Value = 7500
$ws.Range("P389").Value = 500

# Row 390
$ws.Range("D390").Value = 44468
$ws.Range("K390").Value = 33000
$ws.Range("L390").Value = 34000
$ws.Range("M390").Value = 33500
$ws.Range("P390").Value = 2233

# Row 391
$ws.Range("D391").Value = 44468
$ws.Range("J391").Value = 120
$ws.Range("K391").Value = 31000
$ws.Range("L391").Value = 32000
$ws.Range("M391").Value = 31500
$ws.Range("P391").Value = 2100

# Row 392
$ws.Range("D392").Value = 44468
$ws.Range("J392").Value = 120
$ws.Range("K392").Value = 28000
$ws.Range("L392").Value = 29000
$ws.Range("M392").Value = 28500
$ws.Range("P392").Value = 1900

# Row 393
$ws.Range("H393").Value = "Zafiro rojo"
$ws.Range("J393").Value = 120
$ws.Range("K393").Value = 13000
$ws.Range("L393").Value = 14000
$ws.Range("M393").Value = 13500
$ws.Range("P393").Value = 900

# Row 394
$ws.Range("H394").Value = "Zafiro rojo"
$ws.Range("J394").Value = 130
$ws.Range("K394").Value = 11000
$ws.Range("L394").Value = 12000
$ws.Range("M394").Value = 11500
$ws.Range("P394").Value = 767

# Row 395
$ws.Range("H395").Value = "Zafiro rojo"
$ws.Range("K395").Value = 9000
$ws.Range("L395").Value = 10000
$ws.Range("M395").Value = 9500
$ws.Range("P395").Value = 633

# Row 396
$ws.Range("D396").Value = 44389
$ws.Range("H396").Value = "Zafiro verde"
$ws.Range("J396").Value = 100
$ws.Range("K396").Value = 10000
$ws.Range("L396").Value = 11000
$ws.Range("M396").Value = 10500
$ws.Range("P396").Value = 700

# Row 397
$ws.Range("D397").Value = 44389
$ws.Range("H397").Value = "Zafiro verde"
$ws.Range("J397").Value = 120
$ws.Range("K397").Value = 8000
$ws.Range("L397").Value = 9000
$ws.Range("M397").Value = 8500
$ws.Range("P397").Value = 567

# Row 398
$ws.Range("D398").Value = 44389
$ws.Range("H398").Value = "Zafiro verde"
$ws.Range("J398").Value = 140
$ws.Range("K398").Value = 6000
$ws.Range("L398").Value = 7000
$ws.Range("M398").Value = 6500
$ws.Range("P398").Value = 433

# Row 399
$ws.Range("H399").Value = "Zafiro rojo"
$ws.Range("K399").Value = 23000
$ws.Range("L399").Value = 24000
$ws.Range("M399").Value = 23500
$ws.Range("P399").Value = 1567

# Row 400
$ws.Range("H400").Value = "Zafiro rojo"
$ws.Range("K400").Value = 21000
$ws.Range("L400").Value = 22000
$ws.Range("M400").Value = 21500
$ws.Range("P400").Value = 1433

# Row 401
$ws.Range("H401").Value = "Zafiro rojo"
$ws.Range("J401").Value = 160
$ws.Range("K401").Value = 19000
$ws.Range("L401").Value = 20000
$ws.Range("M401").Value = 19500
$ws.Range("P401").Value = 1300

# Row 402
$ws.Range("D402").Value = 44340
$ws.Range("H402").Value = "Zafiro verde"
$ws.Range("K402").Value = 9000
$ws.Range("L402").Value = 10000
$ws.Range("M402").Value = 9500
$ws.Range("P402").Value = 633

# Row 403
$ws.Range("D403").Value = 44340
$ws.Range("H403").Value = "Zafiro verde"
$ws.Range("J403").Value = 160
$ws.Range("K403").Value = 7000
$ws.Range("L403").Value = 8000
$ws.Range("M403").Value = 7500
$ws.Range("P403").Value = 500

# Row 404
$ws.Range("D404").Value = 44340
$ws.Range("H404").Value = "Zafiro verde"
$ws.Range("K404").Value = 6000
$ws.Range("L404").Value = 7000
$ws.Range("M404").Value = 6500
$ws.Range("P404").Value = 433

# Row 405
$ws.Range("H405").Value = "Zafiro rojo"
$ws.Range("K405").Value = 24000
$ws.Range("L405").Value = 25000
$ws.Range("M405").Value = 24500
$ws.Range("P405").Value = 1633

# Row 406
$ws.Range("H406").Value = "Zafiro rojo"
$ws.Range("J406").Value = 140
$ws.Range("K406").Value = 21000
$ws.Range("L406").Value = 22000
$ws.Range("M406").Value = 21500
$ws.Range("P406").Value = 1433

# Row 407
$ws.Range("H407").Value = "Zafiro rojo"
$ws.Range("K407").Value = 17000
$ws.Range("L407").Value = 18000
$ws.Range("M407").Value = 17500
$ws.Range("P407").Value = 1167

# Row 408
$ws.Range("D408").Value = 44330

# Row 409
$ws.Range("D409").Value = 44330
$ws.Range("J409").Value = 120

# Row 410
$ws.Range("D410").Value = 44330
$ws.Range("J410").Value = 120

# Row 411
$ws.Range("D411").Value = 44271
$ws.Range("H411").Value = "Zafiro verde"
$ws.Range("K411").Value = 7000
$ws.Range("L411").Value = 8000
$ws.Range("M411").Value = 7500
$ws.Range("P411").Value = 500

# Row 412
$ws.Range("D412").Value = 44271
$ws.Range("H412").Value = "Zafiro verde"
$ws.Range("J412").Value = 160
$ws.Range("K412").Value = 6000
$ws.Range("L412").Value = 7000
$ws.Range("M412").Value = 6500
$ws.Range("P412").Value = 433

# Row 413
$ws.Range("D413").Value = 44271
$ws.Range("H413").Value = "Zafiro verde"
$ws.Range("K413").Value = 5000
$ws.Range("L413").Value = 6000
$ws.Range("M413").Value = 5500
$ws.Range("P413").Value = 367

# Row 414
$ws.Range("H414").Value = "Zafiro rojo"
$ws.Range("K414").Value = 14000
$ws.Range("L414").Value = 15000
$ws.Range("M414").Value = 14500
$ws.Range("P414").Value = 967

# Row 415
$ws.Range("H415").Value = "Zafiro rojo"
$ws.Range("J415").Value = 140
$ws.Range("K415").Value = 12000
$ws.Range("L415").Value = 13000
$ws.Range("M415").Value = 12500
$ws.Range("P415").Value = 833

# Row 416
$ws.Range("H416").Value = "Zafiro rojo"
$ws.Range("K416").Value = 10000
$ws.Range("L416").Value = 11000
$ws.Range("M416").Value = 10500
$ws.Range("P416").Value = 700

# Row 417
$ws.Range("D417").Value = 44400
$ws.Range("H417").Value = "Zafiro verde"
$ws.Range("K417").Value = 10000
$ws.Range("L417").Value = 11000
$ws.Range("M417").Value = 10500
$ws.Range("P417").Value = 700

# Row 418
$ws.Range("D418").Value = 44400
$ws.Range("H418").Value = "Zafiro verde"
$ws.Range("J418").Value = 120
$ws.Range("K418").Value = 8000
$ws.Range("L418").Value = 9000
$ws.Range("M418").Value = 8500
$ws.Range("P418").Value = 567

# Row 419
$ws.Range("D419").Value = 44400
$ws.Range("H419").Value = "Zafiro verde"
$ws.Range("K419").Value = 6000
$ws.Range("L419").Value = 7000
$ws.Range("M419").Value = 6500
$ws.Range("P419").Value = 433

# Row 420
$ws.Range("H420").Value = "Zafiro rojo"
$ws.Range("K420").Value = 14000
$ws.Range("L420").Value = 15000
$ws.Range("M420").Value = 14500
$ws.Range("P420").Value = 967

# Row 421
$ws.Range("H421").Value = "Zafiro rojo"
$ws.Range("K421").Value = 12000
$ws.Range("L421").Value = 13000
$ws.Range("M421").Value = 12571
$ws.Range("P421").Value = 838

# Row 422
$ws.Range("H422").Value = "Zafiro rojo"
$ws.Range("J422").Value = 160
$ws.Range("K422").Value = 10000
$ws.Range("L422").Value = 11000
$ws.Range("M422").Value = 10500
$ws.Range("P422").Value = 700

# --- Add new rows 423-425 (appended at the end after the shift) ---
# Row 423
$ws.Range("A423").Value = 1
$ws.Range("B423").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C423").Value = "Arica y Parinacota"
$ws.Range("D423").Value = 44309
$ws.Range("E423").Value = 15
$ws.Range("F423").Value = 100112002
$ws.Range("G423").Value = "Pimiento"
$ws.Range("H423").Value = "Zafiro verde"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 120
$ws.Range("K423").Value = 8000
$ws.Range("L423").Value = 9000
$ws.Range("M423").Value = 8500
$ws.Range("N423").Value = "`$/caja 15 kilos"
$ws.Range("O423").Value = "Región de Arica y Parinacota"
$ws.Range("P423").Value = 567
$ws.Range("Q423").Value = 15
$ws.Range("R423").Value = "Hortaliza"
$ws.Range("D423").NumberFormat = $ws.Range("D374").NumberFormat

# Row 424
$ws.Range("A424").Value = 1
$ws.Range("B424").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C424").Value = "Arica y Parinacota"
$ws.Range("D424").Value = 44309
$ws.Range("E424").Value = 15
$ws.Range("F424").Value = 100112002
$ws.Range("G424").Value = "Pimiento"
$ws.Range("H424").Value = "Zafiro verde"
$ws.Range("I424").Value = "Segunda"
$ws.Range("J424").Value = 140
$ws.Range("K424").Value = 7000
$ws.Range("L424").Value = 7000
$ws.Range("M424").Value = 7000
$ws.Range("N424").Value = "`$/caja 15 kilos"
$ws.Range("O424").Value = "Región de Arica y Parinacota"
$ws.Range("P424").Value = 467
$ws.Range("Q424").Value = 15
$ws.Range("R424").Value = "Hortaliza"
$ws.Range("D424").NumberFormat = $ws.Range("D374").NumberFormat

# Row 425
$ws.Range("A425").Value = 1
$ws.Range("B425").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C425").Value = "Arica y Parinacota"
$ws.Range("D425").Value = 44309
$ws.Range("E425").Value = 15
$ws.Range("F425").Value = 100112002
$ws.Range("G425").Value = "Pimiento"
$ws.Range("H425").Value = "Zafiro verde"
$ws.Range("I425").Value = "Tercera"
$ws.Range("J425").Value = 120
$ws.Range("K425").Value = 5000
$ws.Range("L425").Value = 6000
$ws.Range("M425").Value = 5500
$ws.Range("N425").Value = "`$/caja 15 kilos"
$ws.Range("O425").Value = "Región de Arica y Parinacota"
$ws.Range("P425").Value = 367
$ws.Range("Q425").Value = 15
$ws.Range("R425").Value = "Hortaliza"
$ws.Range("D425").NumberFormat = $ws.Range("D374").NumberFormat
